$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Export as TSV")
$listWs = $wb.Worksheets.Item("data_collection_mode list")

# Add the two new allowed values to the data_collection_mode list sheet.
$listWs.Range("A3").Value = "MRM"
$listWs.Range("A4").Value = "PRM"

# Update the data validation on column X (data_collection_mode) to cover the
# extended list and updated error message.
$rng = $ws.Range("X2:X1048576")
$rng.Validation.Modify(3, 1, 1, "='data_collection_mode list'!`$A`$1:`$A`$4")
$rng.Validation.ErrorTitle = "Value must come from list"
$rng.Validation.ErrorMessage = "Value must be one of: DDA / DIA / MRM / PRM."

# Update the explanatory comment on the header cell to mention the new modes.
$comment = $ws.Range("X1").Comment
[void]$comment.Text("Mode of data collection in tandem MS assays. Either DDA (Data-dependent acquisition), DIA (Data-independent acquisition), MRM (multiple reaction monitoring), or PRM (parallel reaction monitoring).")
